# Generate Report for handoff
# Stamp the "Latest Handoff Datetime" (column D) for the most recently
# handed-off file (row 5 -> 45cefd5f-...) on each locale report sheet.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-25 03:07:25"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-25 03:07:35"
